$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3211.111
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H125").Value = 1079.6538
$ws.Range("I125").Value = 707.3684
$ws.Range("J125").Value = 2090.1428
$ws.Range("K125").Value = 6366.3156
$ws.Range("L125").Value = 18811.2852
$ws.Range("M125").Value = -3906.3156
$ws.Range("N125").Value = -23731.2852
$ws.Range("H135").Value = 11780843
$ws.Range("I135").Value = 3337.1035
$ws.Range("J135").Value = 33127572
$ws.Range("K135").Value = 30033.9315
$ws.Range("L135").Value = 298148148
$ws.Range("M135").Value = -27498.9315
$ws.Range("N135").Value = -298153218

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 9508.75
$ws.Range("J37").Value = 10576.571
$ws.Range("L37").Value = 10576.571
$ws.Range("N37").Value = -11122.571
$ws.Range("H61").Value = 18522650
$ws.Range("I61").Value = 27781028
$ws.Range("J61").Value = 5895
$ws.Range("K61").Value = 27781028
$ws.Range("L61").Value = 5895
$ws.Range("M61").Value = -27780816
$ws.Range("N61").Value = -6319
$ws.Range("H63").Value = 3813.7334
$ws.Range("I63").Value = 2033.3334
$ws.Range("J63").Value = 6484.3335
$ws.Range("K63").Value = 2033.3334
$ws.Range("L63").Value = 6484.3335
$ws.Range("M63").Value = -1347.3334
$ws.Range("N63").Value = -7856.3335
$ws.Range("H66").Value = 3813.7334
$ws.Range("I66").Value = 2033.3334
$ws.Range("J66").Value = 6484.3335
$ws.Range("K66").Value = 10166.667
$ws.Range("L66").Value = 32421.6675
$ws.Range("M66").Value = -6734.666999999999
$ws.Range("N66").Value = -39285.6675
$ws.Range("H74").Value = 16135323
$ws.Range("I74").Value = 23810842
$ws.Range("J74").Value = 16732.8
$ws.Range("K74").Value = 23810842
$ws.Range("L74").Value = 16732.8
$ws.Range("M74").Value = -23809968
$ws.Range("N74").Value = -18480.8
$ws.Range("H77").Value = 16135323
$ws.Range("I77").Value = 23810842
$ws.Range("J77").Value = 16732.8
$ws.Range("K77").Value = 119054210
$ws.Range("L77").Value = 83664
$ws.Range("M77").Value = -119049842
$ws.Range("N77").Value = -92400
$ws.Range("H97").Value = 799.21875
$ws.Range("I97").Value = 702.88464
$ws.Range("K97").Value = 702.88464
$ws.Range("M97").Value = -206.88464
$ws.Range("H136").Value = 18522650
$ws.Range("I136").Value = 27781028
$ws.Range("J136").Value = 5895
$ws.Range("K136").Value = 83343084
$ws.Range("L136").Value = 17685
$ws.Range("M136").Value = -83340534
$ws.Range("N136").Value = -22785

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 1825.7778
$ws.Range("I24").Value = 1825.7778
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1825.7778
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -1590.7778
$ws.Range("N24").ClearContents()
$ws.Range("H94").Value = 200.42857
$ws.Range("I94").Value = 200.42857
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 200.42857
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 250.57143
$ws.Range("N94").ClearContents()
$ws.Range("H141").Value = 58616.9
$ws.Range("J141").Value = 53932.5
$ws.Range("L141").Value = 53932.5
$ws.Range("N141").Value = -64292.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2396.4285
$ws.Range("I31").Value = 1197.0358
$ws.Range("J31").Value = 7194
$ws.Range("K31").Value = 1197.0358
$ws.Range("L31").Value = 7194
$ws.Range("M31").Value = -902.0358000000001
$ws.Range("N31").Value = -7784
$ws.Range("H34").Value = 2396.4285
$ws.Range("I34").Value = 1197.0358
$ws.Range("J34").Value = 7194
$ws.Range("K34").Value = 1197.0358
$ws.Range("L34").Value = 7194
$ws.Range("M34").Value = -995.0358000000001
$ws.Range("N34").Value = -7598
$ws.Range("H50").Value = 10596.833
$ws.Range("J50").Value = 10596.833
$ws.Range("L50").Value = 10596.833
$ws.Range("N50").Value = -11846.833
$ws.Range("H51").Value = 9617.182000000001
$ws.Range("J51").Value = 9798.777
$ws.Range("L51").Value = 9798.777
$ws.Range("N51").Value = -11270.777
$ws.Range("H59").Value = 15842.857
$ws.Range("J59").Value = 15816.667
$ws.Range("L59").Value = 15816.667
$ws.Range("N59").Value = -18106.667
$ws.Range("H60").Value = 10201.714
$ws.Range("J60").Value = 10201.714
$ws.Range("L60").Value = 10201.714
$ws.Range("N60").Value = -11223.714
$ws.Range("H61").Value = 9617.182000000001
$ws.Range("J61").Value = 9798.777
$ws.Range("L61").Value = 9798.777
$ws.Range("N61").Value = -10494.777
$ws.Range("H68").Value = 18591.4
$ws.Range("J68").Value = 19672.25
$ws.Range("L68").Value = 19672.25
$ws.Range("N68").Value = -21170.25
$ws.Range("H71").Value = 18591.4
$ws.Range("J71").Value = 19672.25
$ws.Range("L71").Value = 59016.75
$ws.Range("N71").Value = -66504.75
$ws.Range("H74").Value = 14409.223
$ws.Range("J74").Value = 17044.715
$ws.Range("L74").Value = 17044.715
$ws.Range("N74").Value = -18792.715
$ws.Range("H77").Value = 14409.223
$ws.Range("J77").Value = 17044.715
$ws.Range("L77").Value = 51134.145
$ws.Range("N77").Value = -59870.145
$ws.Range("H141").Value = 74887.664
$ws.Range("J141").Value = 74887.664
$ws.Range("L141").Value = 74887.664
$ws.Range("N141").Value = -85247.664

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5024.2383
$ws.Range("I70").Value = 5072.375
$ws.Range("J70").Value = 4994.615
$ws.Range("K70").Value = 5072.375
$ws.Range("L70").Value = 4994.615
$ws.Range("M70").Value = -4802.375
$ws.Range("N70").Value = -5534.615
$ws.Range("H73").Value = 5024.2383
$ws.Range("I73").Value = 5072.375
$ws.Range("J73").Value = 4994.615
$ws.Range("K73").Value = 5072.375
$ws.Range("L73").Value = 4994.615
$ws.Range("M73").Value = -4136.375
$ws.Range("N73").Value = -6866.615
$ws.Range("H131").Value = 18610
$ws.Range("J131").Value = 18610
$ws.Range("L131").Value = 18610
$ws.Range("N131").Value = -28690
$ws.Range("H132").Value = 2552.05
$ws.Range("I132").Value = 2131.9412
$ws.Range("K132").Value = 6395.823600000001
$ws.Range("M132").Value = -3865.823600000001
$ws.Range("H136").Value = 45332.145
$ws.Range("J136").Value = 45332.145
$ws.Range("L136").Value = 135996.435
$ws.Range("N136").Value = -141096.435

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2653.4119
$ws.Range("I40").Value = 2388.5
$ws.Range("J40").Value = 2888.889
$ws.Range("K40").Value = 2388.5
$ws.Range("L40").Value = 2888.889
$ws.Range("M40").Value = -2252.5
$ws.Range("N40").Value = -3160.889
$ws.Range("H68").Value = 2293.6843
$ws.Range("I68").Value = 1972.5
$ws.Range("K68").Value = 1972.5
$ws.Range("M68").Value = -1223.5
$ws.Range("H71").Value = 2293.6843
$ws.Range("I71").Value = 1972.5
$ws.Range("K71").Value = 9862.5
$ws.Range("M71").Value = -6118.5
$ws.Range("H136").Value = 2207.5
$ws.Range("I136").Value = 1345.1111
$ws.Range("J136").Value = 3759.8
$ws.Range("K136").Value = 4035.3333
$ws.Range("L136").Value = 11279.4
$ws.Range("M136").Value = -1485.3333
$ws.Range("N136").Value = -16379.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8242.789000000001
$ws.Range("I136").Value = 12302.869
$ws.Range("J136").Value = 2017.3334
$ws.Range("K136").Value = 36908.607
$ws.Range("L136").Value = 6052.0002
$ws.Range("M136").Value = -34358.607
$ws.Range("N136").Value = -11152.0002
